# Update res_bus vm_pu.xlsx for "case with 380 kV done"
# Columns B-F and I-N for rows 2-25 get new per-unit voltage magnitude values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B,C,D,E,F (rows 2-25)
$valuesBF = @(
    @(1.02,1.03148541099233,1.047306027732941,1.041538956723896,1.053666922211869),
    @(1.02,1.03225628289602,1.047860804491825,1.042208037903043,1.054379481876735),
    @(1.02,1.032755916406945,1.048220379435842,1.042642085963329,1.054841681012223),
    @(1.02,1.032966158551251,1.048371685663215,1.042824822879325,1.055036256677122),
    @(1.02,1.033001470581343,1.048397098847955,1.042855520579532,1.055068942387037),
    @(1.02,1.03275872490424,1.048222400645947,1.042644526671738,1.054844279894932),
    @(1.02,1.031745758165606,1.047493392221064,1.04176484520412,1.053907500703907),
    @(1.02,1.029967208931432,1.046213452194299,1.0402233112041,1.052265495606573),
    @(1.02,1.028785949077922,1.045363419898123,1.03920151426416,1.051176824099277),
    @(1.02,1.028275528578907,1.044996146351774,1.038760489202055,1.050706869014342),
    @(1.02,1.028086098892599,1.044859846355802,1.038596888346527,1.050532526514321),
    @(1.02,1.028126724784893,1.044889077623545,1.038631971498985,1.050569913590564),
    @(1.02,1.02825986692066,1.044984877255443,1.038746961495756,1.050692453312525),
    @(1.02,1.028341921873037,1.04504391876812,1.038817839247662,1.050767983268435),
    @(1.02,1.028819846782056,1.045387811589676,1.039230813744444,1.051208044156506),
    @(1.02,1.029119925109687,1.045603741118494,1.039490243389657,1.051484471764764),
    @(1.02,1.029295059095539,1.045729765822209,1.039641701151904,1.051645846754441),
    @(1.02,1.029354792731013,1.045772749927107,1.039693367478553,1.051700895035373),
    @(1.02,1.029087718837891,1.045580565978459,1.039462394874237,1.051454799258209),
    @(1.02,1.028220655358312,1.044956663273363,1.03871309384166,1.05065636233364),
    @(1.02,1.027676443446794,1.044565096362127,1.038243226376201,1.0501556258854),
    @(1.02,1.027964850193158,1.044772605800836,1.03849219300049,1.05042095433936),
    @(1.02,1.029102271146678,1.045591037590101,1.039474977997122,1.051468206556157),
    @(1.02,1.030426232677218,1.046543781643365,1.040620805302354,1.052688945942984)
)

# New values for columns I,J,K,L,M,N (rows 2-25)
$valuesIN = @(
    @(1.03892164668376,1.036620876777417,1.050069034670119,1.044318186378461,1.056412259242944,1.016195276744317),
    @(1.039060588742105,1.037034553001851,1.05043629598212,1.044798287530243,1.056938166464557,1.01633354031845),
    @(1.03914965874165,1.037302369039882,1.050673827858964,1.045109368745053,1.057278912619181,1.016423020380937),
    @(1.039186903129764,1.037414991012139,1.050773658257039,1.045240247121191,1.057422267987863,1.016460640722366),
    @(1.039193144837967,1.037433902579854,1.050790418534707,1.045262227984362,1.057446344129819,1.016466957500195),
    @(1.039150157192345,1.037303873775741,1.050675161909535,1.045111117158047,1.057280827727641,1.016423523054919),
    @(1.038968775181392,1.036760650810426,1.050193174386736,1.044480350130304,1.056589897991452,1.01624200031041),
    @(1.03864279814639,1.03580455442317,1.049343068427164,1.043372173268251,1.055375910486719,1.015922265402013),
    @(1.03842125205051,1.03516800004771,1.04877589430015,1.042635713068182,1.054569059849928,1.015709226830794),
    @(1.038324327316763,1.03489258241788,1.048530216979875,1.042317389241559,1.054220294535419,1.015617013077928),
    @(1.038288176648532,1.034790313708454,1.048438949927235,1.042199236632147,1.054090840543617,1.015582766307658),
    @(1.038295937796924,1.034812249147791,1.04845852751811,1.042224576803002,1.054118604626847,1.015590112097974),
    @(1.038321342116093,1.034884128159581,1.048522673043693,1.042307620925008,1.054209591920556,1.015614182113579),
    @(1.038336974901044,1.034928419693883,1.048562193733688,1.042358798731505,1.054265664571402,1.015629013191812),
    @(1.038427663755873,1.035186283251477,1.04879219738956,1.04265685128672,1.054592219190535,1.015715347500353),
    @(1.038484284972984,1.03534809270826,1.048936450354003,1.042843965097467,1.054797221936776,1.015769512022922),
    @(1.038517215228626,1.035442494039052,1.049020582152736,1.042953160139559,1.054916854946975,1.015801108456439),
    @(1.038528427296919,1.03547468590635,1.049049267430767,1.042990402058551,1.054957656556329,1.015811882548191),
    @(1.038478219970233,1.035330729942656,1.048920974245739,1.042823883884336,1.054775221030506,1.015763700349239),
    @(1.038313865270343,1.034862960641551,1.048503784091866,1.042283164085098,1.054182795848718,1.015607093937331),
    @(1.038209670132154,1.034569050970888,1.04824141397892,1.041943696449533,1.053810853475333,1.015508661479959),
    @(1.038264987083879,1.034724839048305,1.048380507092144,1.042123606307085,1.054007975461054,1.015560839164078),
    @(1.038480960780467,1.035338575365566,1.048927967259423,1.042832957551308,1.054785162112915,1.015766326384648),
    @(1.038727819710077,1.036051585649804,1.049562924710507,1.043658260840771,1.055689327160208,1.016004905908802)
)

$startRow = 2
$colsBF = @(2,3,4,5,6)      # B,C,D,E,F
$colsIN = @(9,10,11,12,13,14) # I,J,K,L,M,N

for ($i = 0; $i -lt $valuesBF.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValsBF = $valuesBF[$i]
    for ($j = 0; $j -lt $colsBF.Length; $j++) {
        $ws.Cells.Item($rowIndex, $colsBF[$j]).Value = $rowValsBF[$j]
    }
    $rowValsIN = $valuesIN[$i]
    for ($j = 0; $j -lt $colsIN.Length; $j++) {
        $ws.Cells.Item($rowIndex, $colsIN[$j]).Value = $rowValsIN[$j]
    }
}
